$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 32
$ws.Range("B29").Value = "BST all operations"
$ws.Range("H29").Value = "BinarySearchTree"

$ws.Range("A30").Select()
